# Created multi-modal communication observation and added the appropriate
# code systems and value sets.
#
# - Metadata sheet: bump the generation Date and the concept Count.
# - Concepts sheet: append the two newly-defined concepts
#   (multi-modal-communication, hearing).

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$concepts = $wb.Worksheets.Item("Concepts")

# --- Metadata sheet updates -------------------------------------------------
$meta.Range("B8").Value = "2022-03-14T15:55:03-04:00"   # Date
$meta.Range("B23").Value = "5"                           # Count

# --- Concepts sheet: append new rows ----------------------------------------
# Copy row 4's formatting (border/fill/font/alignment) down into the two new
# rows first, so the new cells share the same style as the existing concept
# rows instead of picking up a blank/default style.
$concepts.Range("A4:D4").Copy($concepts.Range("A5:D5"))
$concepts.Range("A4:D4").Copy($concepts.Range("A6:D6"))

# Row 5: multi-modal-communication
$concepts.Range("B5").Value = "multi-modal-communication"
$concepts.Range("C5").Value = "Multi-modal communication"
$concepts.Range("D5").Value = "Category code for multi-modal communication observation"

# Row 6: hearing
$concepts.Range("B6").Value = "hearing"
$concepts.Range("C6").Value = "Hearing"
$concepts.Range("D6").Value = "Category code for hearing observation"
